$d = $word.ActiveDocument

# Insert the six new paragraphs (two empty N1 spacer paragraphs, the three
# SigBlock signature lines, and a trailing empty N1 spacer paragraph) right
# at the very end of the document body, after the existing last paragraph
# and before the sectPr.
$startOfInsert = $d.Content.End
$insertionPoint = $d.Range($startOfInsert, $startOfInsert)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="N1"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="0"/>
    </w:numPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="N1"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="0"/>
    </w:numPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="SigBlock"/>
  </w:pPr>
  <w:r>
    <w:tab/>
  </w:r>
  <w:r>
    <w:t>Shane Doris</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="SigBlock"/>
    <w:rPr>
      <w:rStyle w:val="Sigtitle"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:tab/>
  </w:r>
  <w:r>
    <w:t>A senior officer of the</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="SigBlock"/>
  </w:pPr>
  <w:r>
    <w:tab/>
  </w:r>
  <w:r>
    <w:t>Department of Agriculture, Environment and Rural Affairs</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="N1"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="0"/>
    </w:numPr>
  </w:pPr>
</w:p>
</w:body></w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$insertionPoint.InsertXML($xml)

# The raw-OOXML insertion above does not carry run-level character styles
# (w:r/w:rPr/w:rStyle) through, so (re)apply the three character styles on
# the exact run text now that the paragraphs/text exist. Scope each Find to
# the newly-inserted tail of the document (starting at $startOfInsert) so it
# can't accidentally match similar wording earlier in the body (e.g. "...
# Department of Agriculture, Environment and Rural Affairs..." already
# appears in an earlier paragraph).
$find1 = $d.Range($startOfInsert, $d.Content.End)
$find1.Find.ClearFormatting()
$find1.Find.Execute("Shane Doris", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$find1.Style = "SigSignee"

$find2 = $d.Range($startOfInsert, $d.Content.End)
$find2.Find.ClearFormatting()
$find2.Find.Execute("A senior officer of the", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$find2.Style = "Sigtitle"

$find3 = $d.Range($startOfInsert, $d.Content.End)
$find3.Find.ClearFormatting()
$find3.Find.Execute("Department of Agriculture, Environment and Rural Affairs", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$find3.Style = "Sigtitle"

Write-Output "done"
